{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the start paragraph (\"The objects are stored in HashMaps...\") and\n// the end paragraph (\"...is allowed to leave one review.\") by distinctive\n// text, so the script does not depend on a brittle fixed paragraph index.\nlet startPara = null;\nlet endPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"The objects are stored in HashMaps\") === 0) {\n    startPara = p;\n  }\n  if (p.text.indexOf(\"is allowed to leave one review.\") !== -1) {\n    endPara = p;\n  }\n}\n\nif (!startPara || !endPara) {\n  throw new Error(\"Could not locate the target paragraphs to replace.\");\n}\n\nconst range = startPara.getRange(\"Start\").expandTo(endPara.getRange(\"End\"));\n\nconst flatOpc = `<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w:rsidR=\"006701CF\" w:rsidRDefault=\"006701CF\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">The objects are stored in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>HashMaps</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">. The </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hashmap</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> is stored in the cinema class. If the system was to cater for multiple cinema instances, some of these </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hashmaps</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> would </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>be located in</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> their classes. For example, the same film will be shown in multiple cinemas, so if the system was to cater for multiple cinemas then the collection of films is stored in the film class to allow the user to compare films across multiple cinemas. </w:t></w:r></w:p><w:p w:rsidR=\"0062178D\" w:rsidRDefault=\"0062178D\" w:rsidP=\"00F32A68\"/><w:p w:rsidR=\"0062178D\" w:rsidRDefault=\"0062178D\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">Film stores </w:t></w:r><w:r w:rsidR=\"00A21448\"><w:t>subtitles as empty string if there are no subtitles</w:t></w:r><w:r w:rsidR=\"00C35308\"><w:t xml:space=\"preserve\">. The constructor and methods are overloaded if there are no subtitles for the film. </w:t></w:r></w:p><w:p w:rsidR=\"005E30B9\" w:rsidRDefault=\"005E30B9\" w:rsidP=\"00F32A68\"/><w:p w:rsidR=\"005E30B9\" w:rsidRDefault=\"005E30B9\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">The rating is </w:t></w:r><w:r w:rsidR=\"0042272B\"><w:t xml:space=\"preserve\">stored in the seat-assignment. Since bookings can be for more than one ticket, each ticket </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>is allowed to</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> leave one review. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">Projections -&gt; before a projection is created, it must first be validated to ensure that the film and the screen are available. This is carried out by the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>isValidProjection</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> method. This validation process will still work if there are multiple copies of the film. Each copy of the film is assigned to a separate object. Since the method checks whether the objects are equal, rather than whether the titles of the films are equal, the method will work correctly for multiple copies of the same film. </w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\nrange.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the start paragraph (\"The objects are stored in HashMaps...\") and the\n# end paragraph (\"The rating is ... leave one review.\") by matching their\n# distinctive text, so the script is resilient to exact paragraph-index shifts.\n$startPara = $null\n$endPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"The objects are stored in HashMaps*\") {\n        $startPara = $p\n    }\n    if ($t -like \"*is allowed to leave one review.*\") {\n        $endPara = $p\n    }\n}\n\nif ($startPara -eq $null -or $endPara -eq $null) {\n    throw \"Could not locate the target paragraphs to replace.\"\n}\n\n$range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n\n$xml = '<w:p w:rsidR=\"006701CF\" w:rsidRDefault=\"006701CF\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">The objects are stored in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>HashMaps</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">. The </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hashmap</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> is stored in the cinema class. If the system was to cater for multiple cinema instances, some of these </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>hashmaps</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> would </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>be located in</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> their classes. For example, the same film will be shown in multiple cinemas, so if the system was to cater for multiple cinemas then the collection of films is stored in the film class to allow the user to compare films across multiple cinemas. </w:t></w:r></w:p><w:p w:rsidR=\"0062178D\" w:rsidRDefault=\"0062178D\" w:rsidP=\"00F32A68\"/><w:p w:rsidR=\"0062178D\" w:rsidRDefault=\"0062178D\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">Film stores </w:t></w:r><w:r w:rsidR=\"00A21448\"><w:t>subtitles as empty string if there are no subtitles</w:t></w:r><w:r w:rsidR=\"00C35308\"><w:t xml:space=\"preserve\">. The constructor and methods are overloaded if there are no subtitles for the film. </w:t></w:r></w:p><w:p w:rsidR=\"005E30B9\" w:rsidRDefault=\"005E30B9\" w:rsidP=\"00F32A68\"/><w:p w:rsidR=\"005E30B9\" w:rsidRDefault=\"005E30B9\" w:rsidP=\"00F32A68\"><w:r><w:t xml:space=\"preserve\">The rating is </w:t></w:r><w:r w:rsidR=\"0042272B\"><w:t xml:space=\"preserve\">stored in the seat-assignment. Since bookings can be for more than one ticket, each ticket </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>is allowed to</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> leave one review. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">Projections -&gt; before a projection is created, it must first be validated to ensure that the film and the screen are available. This is carried out by the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>isValidProjection</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> method. This validation process will still work if there are multiple copies of the film. Each copy of the film is assigned to a separate object. Since the method checks whether the objects are equal, rather than whether the titles of the films are equal, the method will work correctly for multiple copies of the same film. </w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n\n$range.InsertXML($xml) | Out-Null\n"}
